# Update the weather observation data: shift rows 2 and 3 to the next
# measurement interval (values sourced from the newer export) and drop the
# old trailing row 4, which is no longer part of the dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was row 2 -> now the first interval of the refreshed export)
$ws.Range("B2").Value2 = 45669.35486111111
$ws.Range("C2").Value2 = 45670.31597222222
$ws.Range("D2").Value2 = -2.8
$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 0.17
$ws.Range("G2").Value2 = 0.3

# Row 3 (was row 3 -> now the second interval of the refreshed export)
$ws.Range("B3").Value2 = 45670.31597222222
$ws.Range("C3").Value2 = 45671.336805555555
$ws.Range("D3").Value2 = -5.4
$ws.Range("E3").Value2 = 3.5
$ws.Range("F3").Value2 = -1.48
$ws.Range("G3").Value2 = -2

# Row 4 no longer exists in the refreshed dataset - remove it entirely so
# the sheet's dimension shrinks back down to A1:G3.
$ws.Rows("4:4").Delete()
